$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.116.86'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.608.96'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.14'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.10'
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.207'
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.648'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.86'
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000305'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.55'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.185.71'
$ws.Range("E14").Value = '  +3.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.14'
$ws.Range("E15").Value = '  +4.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '592.88'
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.21'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.251.51'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.605.70'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.69'
$ws.Range("E22").Value = '  -1.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.17'
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.04'
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.60'
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("E26").Value = '  -0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.82'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.79'
$ws.Range("E30").Value = '  +4.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.31'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0888'
$ws.Range("E35").Value = '  +8.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.927.38'
$ws.Range("E37").Value = '  +3.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '529.97'
$ws.Range("E38").Value = '  +6.37%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.15'
$ws.Range("E40").Value = '  +1.36%  '
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.133'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0455'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.38'
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.86'
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("E48").Value = '  -1.22%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000256'
$ws.Range("E50").Value = '  +4.86%  '
$ws.Range("E51").Value = '  +3.60%  '
